$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 41 (pushes existing rows 41.. down by one)
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row with the new product entry
$ws.Cells.Item(41, 1).Value = "10018154"
$ws.Cells.Item(41, 2).Value = "SARI ROTI TW.CHIP275"
$ws.Cells.Item(41, 3).Value = "BAK02D"
$ws.Cells.Item(41, 4).Value = "5"
$ws.Cells.Item(41, 5).Value = "9"
$ws.Cells.Item(41, 6).Value = "RT,(E-1H)"
